# Apply "repull data, push all data, mean calculation" edit:
# Update column F (dSF) values for rows 3-18 (row indices as they repull from source data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = -4
    5  = -1
    6  = 5
    8  = 5
    9  = -6
    10 = -2
    11 = -2
    12 = -3
    13 = 7
    15 = -6
    16 = 1
    17 = -3
    18 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
